# Überzeitkontrolle_Technik.xlsx edit:
#  - Shift the logged entries on the "mandreoli" sheet by one date slot and
#    add a new fourth entry, updating the Zuschlag (surcharge) values too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mandreoli")

# Columns A (Datum) and D (Zeit) hold values that look numeric/date-like
# ("16.01.2019", "0.50", ...). Force them to be stored as text (as they
# were originally) instead of being auto-converted to dates/numbers.
$ws.Range("A3:A5").NumberFormat = "@"
$ws.Range("D3:D5").NumberFormat = "@"

# Row 3
$ws.Range("A3").Value = "16.01.2019"
$ws.Range("D3").Value = "0.50"
$ws.Range("E3").Value = 1

# Row 4
$ws.Range("A4").Value = "20.01.2019"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = 1.25

# Row 5
$ws.Range("A5").Value = "25.01.2019"
$ws.Range("D5").Value = "2.00"
$ws.Range("E5").Value = 1
